$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2: product name text changed (Galactic Silver -> Sonic Black variant)
$ws.Range("C2").Value = "OnePlus 11R 5G (Sonic Black, 8GB RAM, 128GB Storage)"

# F2: was empty (with a plain style), now holds numeric 0 formatted like the
# hyperlink-styled, vertically-centered cells used elsewhere in the row.
$ws.Range("G2").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("F2").Value = 0

# F3: was a plain number 738, now stored as text "738" (quote-prefixed) so it
# round-trips as a shared string rather than a number.
$ws.Range("F3").Value = "'738"

# Column widths: column F widens to match column G (single merged col span).
$ws.Range("F1").ColumnWidth = 15.25

# Sheet view: drop the frozen/pinned top-left cell and move the active
# selection to F3.
$ws.Range("F3").Select() | Out-Null
